$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1876.8572
$ws.Cells.Item(40, 9).Value = 1574.9
$ws.Cells.Item(40, 10).Value = 2631.75
$ws.Cells.Item(40, 11).Value = 1574.9
$ws.Cells.Item(40, 12).Value = 2631.75
$ws.Cells.Item(40, 13).Value = -1399.9
$ws.Cells.Item(40, 14).Value = -2981.75

$ws.Cells.Item(76, 8).Value = 3060
$ws.Cells.Item(76, 9).Value = 3053.3333
$ws.Cells.Item(76, 10).Value = 3080
$ws.Cells.Item(76, 11).Value = 3053.3333
$ws.Cells.Item(76, 12).Value = 3080
$ws.Cells.Item(76, 13).Value = -2738.3333
$ws.Cells.Item(76, 14).Value = -3710

$ws.Cells.Item(79, 8).Value = 3060
$ws.Cells.Item(79, 9).Value = 3053.3333
$ws.Cells.Item(79, 10).Value = 3080
$ws.Cells.Item(79, 11).Value = 3053.3333
$ws.Cells.Item(79, 12).Value = 3080
$ws.Cells.Item(79, 13).Value = -1961.3333
$ws.Cells.Item(79, 14).Value = -5264

$ws.Cells.Item(138, 8).Value = 5585397
$ws.Cells.Item(138, 9).Value = 1512.9286
$ws.Cells.Item(138, 10).Value = 25128992
$ws.Cells.Item(138, 11).Value = 4538.7858
$ws.Cells.Item(138, 12).Value = 75386976
$ws.Cells.Item(138, 13).Value = 601.2142000000003
$ws.Cells.Item(138, 14).Value = -75397256

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6830.7705
$ws.Cells.Item(32, 9).Value = 3915.46
$ws.Cells.Item(32, 10).Value = 20082.182
$ws.Cells.Item(32, 11).Value = 3915.46
$ws.Cells.Item(32, 12).Value = 20082.182
$ws.Cells.Item(32, 13).Value = -3628.46
$ws.Cells.Item(32, 14).Value = -20656.182

$ws.Cells.Item(132, 8).Value = 2379.2856
$ws.Cells.Item(132, 9).Value = 1708.5
$ws.Cells.Item(132, 10).Value = 3514.4614
$ws.Cells.Item(132, 11).Value = 5125.5
$ws.Cells.Item(132, 12).Value = 10543.3842
$ws.Cells.Item(132, 13).Value = -2595.5
$ws.Cells.Item(132, 14).Value = -15603.3842

$ws.Cells.Item(137, 8).Value = 43000
$ws.Cells.Item(137, 9).Value = 70000
$ws.Cells.Item(137, 10).Value = 37600
$ws.Cells.Item(137, 11).Value = 70000
$ws.Cells.Item(137, 12).Value = 37600
$ws.Cells.Item(137, 13).Value = -64900
$ws.Cells.Item(137, 14).Value = -47800

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 32186.154
$ws.Cells.Item(52, 10).Value = 32186.154
$ws.Cells.Item(52, 12).Value = 32186.154
$ws.Cells.Item(52, 14).Value = -32712.154

$ws.Cells.Item(55, 8).Value = 29135.6
$ws.Cells.Item(55, 10).Value = 29135.6
$ws.Cells.Item(55, 12).Value = 29135.6
$ws.Cells.Item(55, 14).Value = -29681.6

$ws.Cells.Item(57, 8).Value = 40272.727
$ws.Cells.Item(57, 10).Value = 40272.727
$ws.Cells.Item(57, 12).Value = 40272.727
$ws.Cells.Item(57, 14).Value = -41712.727

$ws.Cells.Item(116, 8).Value = 28805.428
$ws.Cells.Item(116, 10).Value = 28805.428
$ws.Cells.Item(116, 12).Value = 28805.428
$ws.Cells.Item(116, 14).Value = -37983.428

$ws.Cells.Item(121, 8).Value = 32186.154
$ws.Cells.Item(121, 10).Value = 32186.154
$ws.Cells.Item(121, 12).Value = 32186.154
$ws.Cells.Item(121, 14).Value = -35680.15399999999

$ws.Cells.Item(136, 8).Value = 40272.727
$ws.Cells.Item(136, 10).Value = 40272.727
$ws.Cells.Item(136, 12).Value = 40272.727
$ws.Cells.Item(136, 14).Value = -50472.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 10800
$ws.Cells.Item(25, 9).Value = 6200
$ws.Cells.Item(25, 11).Value = 6200
$ws.Cells.Item(25, 13).Value = -6026

$ws.Cells.Item(31, 8).Value = 1846.8889
$ws.Cells.Item(31, 9).Value = 1389.3
$ws.Cells.Item(31, 10).Value = 2418.875
$ws.Cells.Item(31, 11).Value = 1389.3
$ws.Cells.Item(31, 12).Value = 2418.875
$ws.Cells.Item(31, 13).Value = -1094.3
$ws.Cells.Item(31, 14).Value = -3008.875

$ws.Cells.Item(34, 8).Value = 1846.8889
$ws.Cells.Item(34, 9).Value = 1389.3
$ws.Cells.Item(34, 10).Value = 2418.875
$ws.Cells.Item(34, 11).Value = 1389.3
$ws.Cells.Item(34, 12).Value = 2418.875
$ws.Cells.Item(34, 13).Value = -1187.3
$ws.Cells.Item(34, 14).Value = -2822.875

$ws.Cells.Item(100, 8).Value = 30995
$ws.Cells.Item(100, 10).Value = 30995
$ws.Cells.Item(100, 12).Value = 30995
$ws.Cells.Item(100, 14).Value = -33159

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1242.5264
$ws.Cells.Item(5, 9).Value = 492.8
$ws.Cells.Item(5, 10).Value = 2075.5557
$ws.Cells.Item(5, 11).Value = 1478.4
$ws.Cells.Item(5, 12).Value = 6226.6671
$ws.Cells.Item(5, 13).Value = -1366.4
$ws.Cells.Item(5, 14).Value = -6450.6671

$ws.Cells.Item(135, 8).Value = 1242.5264
$ws.Cells.Item(135, 9).Value = 492.8
$ws.Cells.Item(135, 10).Value = 2075.5557
$ws.Cells.Item(135, 11).Value = 4435.2
$ws.Cells.Item(135, 12).Value = 18680.0013
$ws.Cells.Item(135, 13).Value = -1900.2
$ws.Cells.Item(135, 14).Value = -23750.0013

$ws.Cells.Item(137, 8).Value = 30760386
$ws.Cells.Item(137, 9).Value = 2704.9167
$ws.Cells.Item(137, 10).Value = 67669600
$ws.Cells.Item(137, 11).Value = 8114.750100000001
$ws.Cells.Item(137, 12).Value = 203008800
$ws.Cells.Item(137, 13).Value = -3014.750100000001
$ws.Cells.Item(137, 14).Value = -203019000

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 6660
$ws.Cells.Item(48, 9).Value = 5000
$ws.Cells.Item(48, 10).Value = 7490
$ws.Cells.Item(48, 11).Value = 5000
$ws.Cells.Item(48, 12).Value = 7490
$ws.Cells.Item(48, 13).Value = -4515
$ws.Cells.Item(48, 14).Value = -8460

$ws.Cells.Item(51, 8).Value = 29200
$ws.Cells.Item(51, 10).Value = 29200
$ws.Cells.Item(51, 12).Value = 29200
$ws.Cells.Item(51, 14).Value = -30218

$ws.Cells.Item(57, 8).Value = 18400
$ws.Cells.Item(57, 10).Value = 18400
$ws.Cells.Item(57, 12).Value = 18400
$ws.Cells.Item(57, 14).Value = -20040

$ws.Cells.Item(113, 8).Value = 1744.579
$ws.Cells.Item(113, 9).Value = 1214.8889
$ws.Cells.Item(113, 10).Value = 2221.3
$ws.Cells.Item(113, 11).Value = 1214.8889
$ws.Cells.Item(113, 12).Value = 2221.3
$ws.Cells.Item(113, 13).Value = 955.1111000000001
$ws.Cells.Item(113, 14).Value = -6561.3

$ws.Cells.Item(132, 8).Value = 4548596
$ws.Cells.Item(132, 9).Value = 3065.8235
$ws.Cells.Item(132, 10).Value = 20003400
$ws.Cells.Item(132, 11).Value = 9197.470499999999
$ws.Cells.Item(132, 12).Value = 60010200
$ws.Cells.Item(132, 13).Value = -6667.470499999999
$ws.Cells.Item(132, 14).Value = -60015260

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 580.12
$ws.Cells.Item(55, 9).Value = 225.28572
$ws.Cells.Item(55, 10).Value = 1031.7273
$ws.Cells.Item(55, 11).Value = 225.28572
$ws.Cells.Item(55, 12).Value = 1031.7273
$ws.Cells.Item(55, 13).Value = -52.28572
$ws.Cells.Item(55, 14).Value = -1377.7273

$ws.Cells.Item(68, 8).Value = 18933.834
$ws.Cells.Item(68, 9).Value = 51350
$ws.Cells.Item(68, 10).Value = 2725.75
$ws.Cells.Item(68, 11).Value = 51350
$ws.Cells.Item(68, 12).Value = 2725.75
$ws.Cells.Item(68, 13).Value = -50601
$ws.Cells.Item(68, 14).Value = -4223.75

$ws.Cells.Item(71, 8).Value = 18933.834
$ws.Cells.Item(71, 9).Value = 51350
$ws.Cells.Item(71, 10).Value = 2725.75
$ws.Cells.Item(71, 11).Value = 256750
$ws.Cells.Item(71, 12).Value = 13628.75
$ws.Cells.Item(71, 13).Value = -253006
$ws.Cells.Item(71, 14).Value = -21116.75

$ws.Cells.Item(122, 8).Value = 3483.6667
$ws.Cells.Item(122, 9).Value = 2004
$ws.Cells.Item(122, 10).Value = 3618.182
$ws.Cells.Item(122, 11).Value = 6012
$ws.Cells.Item(122, 12).Value = 10854.546
$ws.Cells.Item(122, 13).Value = -3562
$ws.Cells.Item(122, 14).Value = -15754.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3701.4062
$ws.Cells.Item(132, 9).Value = 4323.316
$ws.Cells.Item(132, 10).Value = 2792.4614
$ws.Cells.Item(132, 11).Value = 12969.948
$ws.Cells.Item(132, 12).Value = 8377.3842
$ws.Cells.Item(132, 13).Value = -10439.948
$ws.Cells.Item(132, 14).Value = -13437.3842
